$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns after the existing data (AD:AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header styling used by the rest of row 1 (e.g. column A1) by
# copying its format (bold, border, centered) onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the season record for every player row (2-66) with the team's record
for ($row = 2; $row -le 66; $row++) {
    $ws.Cells.Item($row, 30).Value = 77
    $ws.Cells.Item($row, 31).Value = 85
    $ws.Cells.Item($row, 32).Value = 0
}
